# Update values in result_data_RandomForest.xlsx (Sheet1) to reflect
# re-run of the RandomForest imputation algorithm ("Update Name of Algo").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -22.37400000000002
$ws.Range("B4").Value = 5.7159
$ws.Range("A6").Value = -22.44620000000002
$ws.Range("C6").Value = -13.7578
$ws.Range("A7").Value = -20.07949999999997
$ws.Range("C7").Value = -12.52300000000001
$ws.Range("A8").Value = -21.87499999999999
$ws.Range("B8").Value = 6.956400000000001
$ws.Range("C8").Value = -13.35149999999999
$ws.Range("B9").Value = 5.274300000000003
$ws.Range("C10").Value = -13.5648
$ws.Range("B12").Value = 5.026799999999998
$ws.Range("C13").Value = -13.82199999999999
$ws.Range("C14").Value = -14.48589999999999
$ws.Range("A16").Value = -21.5768
$ws.Range("C16").Value = -12.39710000000001
$ws.Range("B17").Value = 5.740499999999995
$ws.Range("B18").Value = 5.986599999999997
$ws.Range("B19").Value = 9.426099999999996
$ws.Range("A20").Value = -22.37970000000002
$ws.Range("B20").Value = 5.027099999999997
$ws.Range("A21").Value = -19.65549999999998
$ws.Range("B26").Value = 4.160900000000004
$ws.Range("A28").Value = -21.93699999999999
$ws.Range("A29").Value = -21.37539999999997
$ws.Range("A30").Value = -22.43720000000001
$ws.Range("C30").Value = -14.1597
$ws.Range("B31").Value = 4.255899999999998
$ws.Range("A32").Value = -21.23560000000002
$ws.Range("C37").Value = -12.8643
$ws.Range("B39").Value = 9.564800000000005
$ws.Range("A40").Value = -21.34699999999997
$ws.Range("B40").Value = 6.031999999999999
$ws.Range("C40").Value = -12.4739
$ws.Range("B41").Value = 9.44169999999999
$ws.Range("B42").Value = 9.326499999999992
$ws.Range("B43").Value = 6.342000000000009
$ws.Range("C44").Value = -13.4949
$ws.Range("A46").Value = -21.6064
$ws.Range("B47").Value = 5.728900000000005
$ws.Range("B48").Value = 5.257700000000004
$ws.Range("A51").Value = -21.76169999999999
$ws.Range("A52").Value = -22.2054
$ws.Range("B54").Value = 4.726599999999999
$ws.Range("A57").Value = -21.89200000000002
$ws.Range("A59").Value = -22.18600000000001
$ws.Range("A62").Value = -22.18380000000002
$ws.Range("B62").Value = 5.200300000000001
$ws.Range("B63").Value = 4.877799999999997
$ws.Range("B64").Value = 5.4262
$ws.Range("A66").Value = -21.46550000000001
$ws.Range("C70").Value = -12.1941
$ws.Range("A73").Value = -20.47929999999998
$ws.Range("A74").Value = -21.57989999999998
$ws.Range("B76").Value = 6.254100000000001
$ws.Range("A77").Value = -19.82759999999998
$ws.Range("B81").Value = 5.514800000000005
$ws.Range("B84").Value = 5.839600000000003
$ws.Range("B89").Value = 5.388399999999998
$ws.Range("C89").Value = -14.01489999999999
$ws.Range("C91").Value = -12.67010000000001
$ws.Range("A92").Value = -21.57320000000001
$ws.Range("C93").Value = -10.44499999999999
$ws.Range("B94").Value = 4.699599999999993
$ws.Range("C98").Value = -12.4092
$ws.Range("A100").Value = -22.04940000000002
